$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns B and G to accommodate new study/travel/work labels
# (values tuned so the engine's internal width-snapping reproduces the
# target stored widths of 15 and 16.42578125 as closely as possible)
$ws.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws.Columns.Item(7).ColumnWidth = 15.666666666666666

# Set cell values in the same order the new labels were first typed in,
# so new shared-string entries are appended in the matching order.

# 1) EST.PDS
$ws.Cells.Item(5, 2).Value = "EST.PDS"
$ws.Cells.Item(6, 2).Value = "EST.PDS"
$ws.Cells.Item(18, 5).Value = "EST.PDS"
$ws.Cells.Item(19, 5).Value = "EST.PDS"

# 2) EST.INFOII
$ws.Cells.Item(9, 2).Value = "EST.INFOII"
$ws.Cells.Item(10, 2).Value = "EST.INFOII"

# 3) EST.ROBOTICA
$ws.Cells.Item(16, 2).Value = "EST.ROBOTICA"
$ws.Cells.Item(17, 2).Value = "EST.ROBOTICA"
$ws.Cells.Item(16, 6).Value = "EST.ROBOTICA"
$ws.Cells.Item(17, 6).Value = "EST.ROBOTICA"

# 4) EST INFOII
$ws.Cells.Item(16, 4).Value = "EST INFOII"
$ws.Cells.Item(17, 4).Value = "EST INFOII"

# 5) EST INGLES
$ws.Cells.Item(6, 8).Value = "EST INGLES"

# 6) TRABAJO
$ws.Cells.Item(11, 7).Value = "TRABAJO"
$ws.Cells.Item(11, 8).Value = "TRABAJO"
$ws.Cells.Item(12, 7).Value = "TRABAJO"
$ws.Cells.Item(12, 8).Value = "TRABAJO"
$ws.Cells.Item(13, 7).Value = "TRABAJO"
$ws.Cells.Item(13, 8).Value = "TRABAJO"
$ws.Cells.Item(14, 7).Value = "TRABAJO"
$ws.Cells.Item(14, 8).Value = "TRABAJO"
$ws.Cells.Item(15, 7).Value = "TRABAJO"
$ws.Cells.Item(15, 8).Value = "TRABAJO"
$ws.Cells.Item(16, 7).Value = "TRABAJO"
$ws.Cells.Item(16, 8).Value = "TRABAJO"
$ws.Cells.Item(17, 7).Value = "TRABAJO"
$ws.Cells.Item(17, 8).Value = "TRABAJO"
$ws.Cells.Item(18, 7).Value = "TRABAJO"
$ws.Cells.Item(18, 8).Value = "TRABAJO"
$ws.Cells.Item(19, 7).Value = "TRABAJO"
$ws.Cells.Item(19, 8).Value = "TRABAJO"

# 7) EST. InGLES
$ws.Cells.Item(7, 8).Value = "EST. InGLES"

# 8) VIAJE AL TRABAJO
$ws.Cells.Item(9, 7).Value = "VIAJE AL TRABAJO"
$ws.Cells.Item(10, 7).Value = "VIAJE AL TRABAJO"

# Reused existing strings (already present in the shared string table)
$ws.Cells.Item(5, 7).Value = "VIAJE A LA U"
$ws.Cells.Item(6, 7).Value = "VIAJE A LA U"
$ws.Cells.Item(7, 7).Value = "LAB INTEGRADO"
$ws.Cells.Item(8, 7).Value = "LAB INTEGRADO"
$ws.Cells.Item(16, 5).Value = "VIAJE A CASA"

# Update the active selection to L13 as in the diff
$ws.Range("L13").Select()
